# Update forecast-error table values (ifo GDP component analysis preprocessing).
# Recomputed ME / MAE / MSE / RMSE / SE / N figures for rows Q0..Q8 (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1374531038001369
$ws.Range("C2").Value = 0.9463274094814685
$ws.Range("D2").Value = 4.340354840748894
$ws.Range("E2").Value = 2.083351828364306
$ws.Range("F2").Value = 2.099094096015294
$ws.Range("G2").Value = 52

$ws.Range("B3").Value = 0.1201387303292092
$ws.Range("C3").Value = 0.9941924500794057
$ws.Range("D3").Value = 4.368983621713564
$ws.Range("E3").Value = 2.090211382064877
$ws.Range("F3").Value = 2.10752018100236
$ws.Range("G3").Value = 51

$ws.Range("B4").Value = 0.1570212086422995
$ws.Range("C4").Value = 0.9129929970791744
$ws.Range("D4").Value = 4.115459563428113
$ws.Range("E4").Value = 2.028659548428004
$ws.Range("F4").Value = 2.043107852613974
$ws.Range("G4").Value = 50

$ws.Range("B5").Value = 0.1668826930103534
$ws.Range("C5").Value = 1.018518659023897
$ws.Range("D5").Value = 4.584217853051491
$ws.Range("E5").Value = 2.141078665778418
$ws.Range("F5").Value = 2.156685501776154
$ws.Range("G5").Value = 49

$ws.Range("B6").Value = 0.1498087523549735
$ws.Range("C6").Value = 0.9654930794235904
$ws.Range("D6").Value = 4.38482909777947
$ws.Range("E6").Value = 2.093998351904669
$ws.Range("F6").Value = 2.110735219814999
$ws.Range("G6").Value = 48

$ws.Range("B7").Value = 0.1482159292065307
$ws.Range("C7").Value = 1.037301259063642
$ws.Range("D7").Value = 5.28735183827951
$ws.Range("E7").Value = 2.299424240604484
$ws.Range("F7").Value = 2.327192174260785
$ws.Range("G7").Value = 36

$ws.Range("B8").Value = 0.1958133101793198
$ws.Range("C8").Value = 1.063338432913966
$ws.Range("D8").Value = 5.471053399423887
$ws.Range("E8").Value = 2.339028302399073
$ws.Range("F8").Value = 2.364845904346298
$ws.Range("G8").Value = 35

$ws.Range("B9").Value = 0.1080685240234339
$ws.Range("C9").Value = 1.516053278149146
$ws.Range("D9").Value = 9.891401411246379
$ws.Range("E9").Value = 3.145059842236134
$ws.Range("F9").Value = 3.234328795688273
$ws.Range("G9").Value = 18

$ws.Range("B10").Value = -0.8008619424696254
$ws.Range("C10").Value = 1.150395842953777
$ws.Range("D10").Value = 6.142350263041148
$ws.Range("E10").Value = 2.478376537784594
$ws.Range("F10").Value = 2.459891756431453
$ws.Range("G10").Value = 11
